$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order country labels (rank swaps in the source table) ---
# Georgia overtakes Mozambique (row 108/109 pair)
$ws.Range("A108").Value = "Georgia"
$ws.Range("A109").Value = "Mozambique"

# Santa Lucia overtakes Nueva Caledonia (row 207/208 pair)
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Nueva Caledonia"

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 09:48"

# --- Update statistic rows ---

# Row 7: Rusia
$ws.Range("B7").Value = 1248619
$ws.Range("C7").Value = 11115
$ws.Range("D7").Value = 995275
$ws.Range("E7").Value = 231479
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 202
$ws.Range("H7").Value = 21865

# Row 60: Singapur
$ws.Range("B60").Value = 57840
$ws.Range("C60").Value = 10
$ws.Range("D60").Value = 57612
$ws.Range("E60").Value = 201
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 27

# Row 63: Armenia
$ws.Range("B63").Value = 53755
$ws.Range("C63").Value = 672
$ws.Range("D63").Value = 45110
$ws.Range("E63").Value = 7650
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 995

# Row 73: Afganistan
$ws.Range("B73").Value = 39548
$ws.Range("C73").Value = 62
$ws.Range("D73").Value = 33045
$ws.Range("E73").Value = 5034
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 1469

# Row 77: Hungria
$ws.Range("B77").Value = 33114
$ws.Range("C77").Value = 816
$ws.Range("D77").Value = 9149
$ws.Range("E77").Value = 23088
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 24
$ws.Range("H77").Value = 877

# Row 81: Australia
$ws.Range("B81").Value = 27182
$ws.Range("C81").Value = 8
$ws.Range("D81").Value = 24937
$ws.Range("E81").Value = 1348
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 897

# Row 108: now Georgia (new, higher data)
$ws.Range("B108").Value = 9753
$ws.Range("C108").Value = 508
$ws.Range("D108").Value = 5235
$ws.Range("E108").Value = 4456
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 4
$ws.Range("H108").Value = 62

# Row 109: now Mozambique (previous Mozambique data, unchanged values, moved down)
$ws.Range("B109").Value = 9398
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 6358
$ws.Range("E109").Value = 2973
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 67
